$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(41).Insert()

$ws.Cells.Item(41, 1).Value = 5
$ws.Cells.Item(41, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(41, 3).Value = "Maule"
$ws.Cells.Item(41, 4).Value = 44414
$ws.Cells.Item(41, 5).Value = 7
$ws.Cells.Item(41, 6).Value = 100114013
$ws.Cells.Item(41, 7).Value = "Zanahoria"
$ws.Cells.Item(41, 8).Value = "Sin especificar"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 400
$ws.Cells.Item(41, 11).Value = 5500
$ws.Cells.Item(41, 12).Value = 5500
$ws.Cells.Item(41, 13).Value = 5500
$ws.Cells.Item(41, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(41, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(41, 16).Value = 275
$ws.Cells.Item(41, 17).Value = 20
$ws.Cells.Item(41, 18).Value = "Hortaliza"
